# Update cryptos list snapshot values (prices + 1h volume %) per Sep 23 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.671.87"
$ws.Range("E2").Value = "  +0.10%  "

# Row 3
$ws.Range("D3").Value = "'1.597.70"
$ws.Range("E3").Value = "  -0.03%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "'211.62"
$ws.Range("E5").Value = "  +0.19%  "

# Row 6
$ws.Range("E6").Value = "  +0.36%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("E8").Value = "  +0.25%  "

# Row 9
$ws.Range("E9").Value = "  +0.66%  "

# Row 10
$ws.Range("D10").Value = "'19.52"
$ws.Range("E10").Value = "  -0.58%  "

# Row 11
$ws.Range("E11").Value = "  +0.21%  "

# Row 12
$ws.Range("D12").Value = "'1.823.58"
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$ws.Range("D13").Value = "'1.576.17"
$ws.Range("E13").Value = "  -2.15%  "

# Row 14
$ws.Range("E14").Value = "  +0.04%  "

# Row 15
$ws.Range("D15").Value = "'0.523"

# Row 16
$ws.Range("D16").Value = "'65.05"
$ws.Range("E16").Value = "  +0.32%  "

# Row 17
$ws.Range("D17").Value = "'26.651.67"
$ws.Range("E17").Value = "  +0.06%  "

# Row 18
$ws.Range("D18").Value = "'0.0₃0738"
$ws.Range("E18").Value = "  +1.38%  "

# Row 19
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'1.00"
$ws.Range("E19").Value = "  +0.19%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'209.08"
$ws.Range("E20").Value = "  +0.05%  "

# Row 21
$ws.Range("E21").Value = "  +3.97%  "

# Row 22
$ws.Range("E22").Value = "  +0.50%  "

# Row 23
$ws.Range("E23").Value = "  +2.97%  "

# Row 24
$ws.Range("E24").Value = "  +1.17%  "

# Row 25
$ws.Range("D25").Value = "'144.16"
$ws.Range("E25").Value = "  -1.25%  "

# Row 26
$ws.Range("E26").Value = "  +0.03%  "

# Row 27
$ws.Range("E27").Value = "  -1.61%  "

# Row 28
$ws.Range("E28").Value = "  -0.78%  "

# Row 29
$ws.Range("D29").Value = "'15.29"
$ws.Range("E29").Value = "  -0.04%  "

# Row 30
$ws.Range("E30").Value = "  +1.75%  "

# Row 31
$ws.Range("E31").Value = "  +0.15%  "

# Row 33
$ws.Range("E33").Value = "  +0.95%  "

# Row 34
$ws.Range("D34").Value = "'1.286.58"
$ws.Range("E34").Value = "  -0.54%  "

# Row 35
$ws.Range("E35").Value = "  -6.84%  "

# Row 36
$ws.Range("E36").Value = "  +0.52%  "

# Row 37
$ws.Range("E37").Value = "  +0.36%  "

# Row 38
$ws.Range("D38").Value = "'0.0171"
$ws.Range("E38").Value = "  -0.15%  "

# Row 39
$ws.Range("D39").Value = "'0.833"
$ws.Range("E39").Value = "  -1.17%  "

# Row 40
$ws.Range("E40").Value = "  +16.04%  "

# Row 41
$ws.Range("D41").Value = "'5.47"
$ws.Range("E41").Value = "  +1.61%  "

# Row 42
$ws.Range("E42").Value = "  -0.49%  "

# Row 43
$ws.Range("D43").Value = "'0.783"
$ws.Range("E43").Value = "  -0.25%  "

# Row 44
$ws.Range("D44").Value = "'63.46"
$ws.Range("E44").Value = "  -0.66%  "

# Row 45
$ws.Range("D45").Value = "'1.733.75"
$ws.Range("E45").Value = "  -0.07%  "

# Row 46
$ws.Range("D46").Value = "'91.07"
$ws.Range("E46").Value = "  +1.35%  "

# Row 47
$ws.Range("E47").Value = "  -2.76%  "

# Row 48
$ws.Range("E48").Value = "  -0.96%  "

# Row 49
$ws.Range("E49").Value = "  +1.48%  "

# Row 50
$ws.Range("E50").Value = "  +0.88%  "

# Row 51
$ws.Range("E51").Value = "  +0.15%  "

